# Updated custom app bar button: refresh scraped product rows with new
# prices and append two newly scraped products (rows 5 and 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Price refreshes on existing rows ----
$ws.Range("F2").Value = 629000
$ws.Range("F4").Value = 549000

# Stash a clean copy of the "Hyperlink" cell format (style shared by
# E2:E4) on a scratch cell well outside the table so we can re-apply it
# later without Excel minting a brand-new style entry.
$ws.Range("E4").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# ---- New row 5: RAM KINGSTON HYPERX FURY ----
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "f8c6fc0f-c983-4edd-8fa7-72454eedcb2c"
$ws.Range("C5").Value = "RAM KINGSTON HYPERX FURY DDR4 8GB 2666MHz 21300 GAMING RAM PC DDR4 8GB"
$ws.Range("D5").Value = "Tokopedia"
$ws.Range("E5").Value = "Tokopedia"
$ws.Range("F5").Value = 390000

# ---- New row 6: RAM DDR4 V-GeN RESCUE ----
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "76ab1cfe-bffc-46e5-b6e1-9f3ecdc24b52"
$ws.Range("C6").Value = "RAM DDR4 V-GeN RESCUE 8GB PC19200/2400Mhz Long Dimm (Memory PC VGEN)"
$ws.Range("D6").Value = "Tokopedia"
$ws.Range("E6").Value = "Tokopedia"
$ws.Range("F6").Value = 458000

# Extend the number/index-column formatting (bold, border, centered) from
# the existing row 4 down onto the two new rows in column A.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A6").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Add the real hyperlinks (shop_link) for the two new rows.
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.tokopedia.com/rajaramnusantara/ram-kingston-hyperx-fury-ddr4-8gb-2666mhz-21300-gaming-ram-pc-ddr4-8gb?src=topads")
$ws.Hyperlinks.Add($ws.Range("E6"), "https://www.tokopedia.com/intact-official/ram-ddr4-v-gen-rescue-8gb-pc19200-2400mhz-long-dimm-memory-pc-vgen?extParam=ivf%3Dfalse%26src%3Dsearch")

# Restore the original "Hyperlink" style (same one used by E2:E4) on the
# new cells, undoing the ad-hoc style Excel just minted for them.
$ws.Range("Z1").Copy()
$ws.Range("E5").PasteSpecial(-4122)
$ws.Range("E6").PasteSpecial(-4122)
$wb.Application.CutCopyMode = $false

# Drop the scratch cell entirely (value + formatting).
$ws.Range("Z1").Clear()
